$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI values for F2 (ligand) -> F2rl1 (receptor) LR pairs
# Values recomputed from refreshed per-cluster average expression levels.

# Row 2
$ws.Range("G2").Value = 0.754521
$ws.Range("H2").Value = 2.263563
$ws.Range("I2").Value = 0.2768403531129761
$ws.Range("J2").Value = 0.2768403531129761
$ws.Range("M2").Value = 0.3423083333333334
$ws.Range("N2").Value = 1.026925
$ws.Range("O2").Value = 0.056115687851804
$ws.Range("P2").Value = 0.056115687851804
$ws.Range("Q2").Value = 0.258278825975
$ws.Range("R2").Value = 2.324509433775
$ws.Range("S2").Value = 0.01553508684007096
$ws.Range("T2").Value = 0.01553508684007096

# Row 3
$ws.Range("G3").Value = 0.754521
$ws.Range("H3").Value = 2.263563
$ws.Range("I3").Value = 0.2768403531129761
$ws.Range("J3").Value = 0.2768403531129761
$ws.Range("O3").Value = 0.943884312148196
$ws.Range("P3").Value = 0.9438843121481961
$ws.Range("Q3").Value = 4.344334736511999
$ws.Range("R3").Value = 39.09901262860799
$ws.Range("S3").Value = 0.2613052662729052
$ws.Range("T3").Value = 0.2613052662729052

# Row 4
$ws.Range("G4").Value = 0.9731926666666667
$ws.Range("I4").Value = 0.3570728998754956
$ws.Range("J4").Value = 0.3570728998754956
$ws.Range("M4").Value = 0.3423083333333334
$ws.Range("N4").Value = 1.026925
$ws.Range("O4").Value = 0.056115687851804
$ws.Range("P4").Value = 0.056115687851804
$ws.Range("Q4").Value = 0.3331319597388889
$ws.Range("R4").Value = 2.99818763765
$ws.Range("S4").Value = 0.02003739138975177
$ws.Range("T4").Value = 0.02003739138975177

# Row 5
$ws.Range("G5").Value = 0.9731926666666667
$ws.Range("I5").Value = 0.3570728998754956
$ws.Range("J5").Value = 0.3570728998754956
$ws.Range("O5").Value = 0.943884312148196
$ws.Range("P5").Value = 0.9438843121481961
$ws.Range("S5").Value = 0.3370355084857439
$ws.Range("T5").Value = 0.3370355084857439

# Row 6
$ws.Range("G6").Value = 0.782441
$ws.Range("H6").Value = 2.347323
$ws.Range("I6").Value = 0.2870844452706686
$ws.Range("J6").Value = 0.2870844452706686
$ws.Range("M6").Value = 0.3423083333333334
$ws.Range("N6").Value = 1.026925
$ws.Range("O6").Value = 0.056115687851804
$ws.Range("P6").Value = 0.056115687851804
$ws.Range("Q6").Value = 0.2678360746416667
$ws.Range("R6").Value = 2.410524671775001
$ws.Range("S6").Value = 0.01610994111791715
$ws.Range("T6").Value = 0.01610994111791715

# Row 7
$ws.Range("G7").Value = 0.782441
$ws.Range("H7").Value = 2.347323
$ws.Range("I7").Value = 0.2870844452706686
$ws.Range("J7").Value = 0.2870844452706686
$ws.Range("O7").Value = 0.943884312148196
$ws.Range("P7").Value = 0.9438843121481961
$ws.Range("Q7").Value = 4.505090800085333
$ws.Range("R7").Value = 40.545817200768
$ws.Range("S7").Value = 0.2709745041527514
$ws.Range("T7").Value = 0.2709745041527515

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2153186666666667
$ws.Range("H8").Value = 0.645956
$ws.Range("I8").Value = 0.07900230174085969
$ws.Range("J8").Value = 0.07900230174085969
$ws.Range("M8").Value = 0.3423083333333334
$ws.Range("N8").Value = 1.026925
$ws.Range("O8").Value = 0.056115687851804
$ws.Range("P8").Value = 0.056115687851804
$ws.Range("Q8").Value = 0.07370537392222223
$ws.Range("R8").Value = 0.6633483653000001
$ws.Range("S8").Value = 0.004433268504064114
$ws.Range("T8").Value = 0.004433268504064114

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2153186666666667
$ws.Range("H9").Value = 0.645956
$ws.Range("I9").Value = 0.07900230174085969
$ws.Range("J9").Value = 0.07900230174085969
$ws.Range("O9").Value = 0.943884312148196
$ws.Range("P9").Value = 0.9438843121481961
$ws.Range("Q9").Value = 1.239748612721778
$ws.Range("R9").Value = 11.157737514496
$ws.Range("S9").Value = 0.07456903323679558
$ws.Range("T9").Value = 0.07456903323679559
